# Apply crypto price/volume updates for Thu May 30 02:56:34 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price (column D) values parse as plain numbers (e.g. "596.41"); the
# source data stores these as literal text (inline strings), matching the other
# Price cells like "68.052.40" that use "." as a thousands separator and cannot
# be parsed as numbers. Force the cell format to Text first so Excel keeps the
# typed characters verbatim instead of auto-converting them to a Number.
$textCells = @("D5","D6","D11","D13","D17","D19","D21","D22","D25","D27","D28","D32","D33","D35","D37","D40","D41","D43","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '67.996.88'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '3.778.04'
$ws.Range("E3").Value = '  -1.85%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '596.41'
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").Value = '170.14'
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("D7").Value = '3.776.29'
$ws.Range("E7").Value = '  -1.99%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("E10").Value = '  -1.66%  '
$ws.Range("D11").Value = '6.52'
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("D13").Value = '0.0000284'
$ws.Range("E13").Value = '  +1.85%  '
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").Value = '4.412.81'
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("D16").Value = '3.749.84'
$ws.Range("E16").Value = '  -1.11%  '
$ws.Range("D17").Value = '18.95'
$ws.Range("E17").Value = '  +3.12%  '
$ws.Range("D18").Value = '67.956.17'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = '7.20'
$ws.Range("E19").Value = '  -2.15%  '
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("D21").Value = '10.61'
$ws.Range("E21").Value = '  -3.07%  '
$ws.Range("D22").Value = '468.04'
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("E23").Value = '  -1.02%  '
$ws.Range("E24").Value = '  -7.38%  '
$ws.Range("D25").Value = '83.71'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '10.56'
$ws.Range("E27").Value = '  +1.17%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '12.15'
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  -1.17%  '
$ws.Range("D31").Value = '3.928.04'
$ws.Range("E31").Value = '  -1.76%  '
$ws.Range("D32").Value = '7.60'
$ws.Range("E32").Value = '  -1.79%  '
$ws.Range("D33").Value = '30.49'
$ws.Range("E33").Value = '  -2.80%  '
$ws.Range("E34").Value = '  -2.44%  '
$ws.Range("D35").Value = '9.21'
$ws.Range("E35").Value = '  -0.76%  '
$ws.Range("D36").Value = '3.738.04'
$ws.Range("E36").Value = '  -1.98%  '
$ws.Range("D37").Value = '3.76'
$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("E38").Value = '  +0.29%  '
$ws.Range("E39").Value = '  -0.35%  '
$ws.Range("D40").Value = '1.01'
$ws.Range("E40").Value = '  -1.51%  '
$ws.Range("D41").Value = '5.84'
$ws.Range("E41").Value = '  -1.45%  '
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").Value = '0.314'
$ws.Range("E43").Value = '  -1.02%  '
$ws.Range("E45").Value = '  -2.31%  '
$ws.Range("D46").Value = '8.68'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").Value = '404.91'
$ws.Range("E47").Value = '  -3.93%  '
$ws.Range("D48").Value = '0.000283'
$ws.Range("E48").Value = '  -4.94%  '
$ws.Range("D49").Value = '45.68'
$ws.Range("E49").Value = '  -2.51%  '
$ws.Range("D50").Value = '40.24'
$ws.Range("E50").Value = '  +7.17%  '
$ws.Range("D51").Value = '140.87'
$ws.Range("E51").Value = '  -0.73%  '
